# leetcode-tracker.xlsx — "String challenge commit 1"
#
# - Rename "Sheet3" -> "String" and populate it with the String-category
#   leetcode problems (title + link per row), mirroring the "Array" sheet's
#   layout (S No / Completed headers in row 1, col widths for B/C/D).
# - Switch the active tab from "Array" to the new "String" sheet, with the
#   selection left on D8.
# - Clear the previously-active selection on "Array" down to a plain header
#   range selection (A1:D1), no longer the active tab.

$wb = $excel.ActiveWorkbook

$arraySheet = $wb.Worksheets.Item("Array")
$stringSheet = $wb.Worksheets.Item("Sheet3")

# 1. Rename Sheet3 -> String
$stringSheet.Name = "String"

# 2. Populate the String sheet — same header layout as Array (S No / Completed)
$stringSheet.Range("A1").Value = "S No"
$stringSheet.Range("D1").Value = "Completed"

$problems = @(
    @("Longest Substring Without Repeating Characters", " - https://leetcode.com/problems/longest-substring-without-repeating-characters/"),
    @("Longest Repeating Character Replacement", " - https://leetcode.com/problems/longest-repeating-character-replacement/"),
    @("Minimum Window Substring", " - https://leetcode.com/problems/minimum-window-substring/"),
    @("Valid Anagram", " - https://leetcode.com/problems/valid-anagram/"),
    @("Group Anagrams", " - https://leetcode.com/problems/group-anagrams/"),
    @("Valid Parentheses", " - https://leetcode.com/problems/valid-parentheses/"),
    @("Valid Palindrome", " - https://leetcode.com/problems/valid-palindrome/"),
    @("Longest Palindromic Substring", " - https://leetcode.com/problems/longest-palindromic-substring/"),
    @("Palindromic Substrings", " - https://leetcode.com/problems/palindromic-substrings/"),
    @("Encode and Decode Strings (Leetcode Premium)", " - https://leetcode.com/problems/encode-and-decode-strings/")
)

$row = 2
foreach ($p in $problems) {
    $stringSheet.Cells.Item($row, 2).Value = $p[0]
    $stringSheet.Cells.Item($row, 3).Value = $p[1]
    $row++
}

# 3. Column widths to match the Array sheet's best-fit look (closest the
#    host's column-width rounding can land to the original 40.77734375 /
#    69.77734375 / 9.88671875 best-fit values)
$stringSheet.Columns.Item(2).ColumnWidth = 40
$stringSheet.Columns.Item(3).ColumnWidth = 69
$stringSheet.Columns.Item(4).ColumnWidth = 9

# 4. Re-point the active tab/selection: String becomes active (D8 selected),
#    Array drops tabSelected and its selection becomes the header row range.
$arraySheet.Range("A1:D1").Select() | Out-Null
$stringSheet.Activate() | Out-Null
$stringSheet.Range("D8").Select() | Out-Null
